$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the serial number text in V2 (shared string used by that cell)
$ws.Range("V2").Value = "38437D095778"

# Update the numeric value in P2
$ws.Range("P2").Value = 509514324

# Update the selection / active cell shown in the sheet view
$ws.Range("V7").Select()
